$d = $word.ActiveDocument

# Update the date/weekday heading paragraph.
$d.Paragraphs.Item(1).Range.Text = "2024-01-11 Thursday"

# Update the 5x5 grid of division problems (only rows 1,5,9,13,17 carry text).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "35÷4="
$t.Cell(1, 2).Range.Text  = "47÷7="
$t.Cell(1, 3).Range.Text  = "40÷2="
$t.Cell(1, 4).Range.Text  = "74÷7="
$t.Cell(1, 5).Range.Text  = "27÷6="

$t.Cell(5, 1).Range.Text  = "88÷8="
$t.Cell(5, 2).Range.Text  = "78÷8="
$t.Cell(5, 3).Range.Text  = "21÷4="
$t.Cell(5, 4).Range.Text  = "80÷3="
$t.Cell(5, 5).Range.Text  = "72÷6="

$t.Cell(9, 1).Range.Text  = "56÷4="
$t.Cell(9, 2).Range.Text  = "98÷7="
$t.Cell(9, 3).Range.Text  = "37÷9="
$t.Cell(9, 4).Range.Text  = "48÷4="
$t.Cell(9, 5).Range.Text  = "26÷8="

$t.Cell(13, 1).Range.Text = "17÷5="
$t.Cell(13, 2).Range.Text = "58÷3="
$t.Cell(13, 3).Range.Text = "90÷3="
$t.Cell(13, 4).Range.Text = "30÷5="
$t.Cell(13, 5).Range.Text = "44÷8="

$t.Cell(17, 1).Range.Text = "22÷5="
$t.Cell(17, 2).Range.Text = "24÷5="
$t.Cell(17, 3).Range.Text = "92÷5="
$t.Cell(17, 4).Range.Text = "33÷4="
$t.Cell(17, 5).Range.Text = "53÷4="
